# Apply update: change sequence number 1000221 -> 1000245 in the
# "Pre-Alert Template Import" sheet, row 3 (test data row).
#
# Cells A3/B3/C3/AN3/AO3 carry a quote-prefix ("Text" stored with a leading
# apostrophe, style index 15 in styles.xml) in the original workbook, so a
# leading "'" is used when re-writing them to keep that formatting intact
# (Excel strips the apostrophe itself; it only flips the quote-prefix flag).
# AJ3/AL3 use a plain (non-quote-prefixed) text style and are written as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pre-Alert Template Import")

$ws.Range("A3").Value = "'JSSO1000245"
$ws.Range("B3").Value = "'JSSO1000245"
$ws.Range("C3").Value = "'JSSO1000245"
$ws.Range("AJ3").Value = "JSCN1000245"
$ws.Range("AL3").Value = "SLJSSO1000245"
$ws.Range("AN3").Value = "'MBLJSSO1000245"
$ws.Range("AO3").Value = "'HBLJSSO1000245"
